$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format first so numeric-looking values (e.g. "0.999")
# are stored as text, matching the original inlineStr cell type.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '46.263.73'
$ws.Range("E2").Value = '  +4.11%  '

# Row 3
$ws.Range("D3").Value = '2.450.97'
$ws.Range("E3").Value = '  +1.37%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '322.54'
$ws.Range("E5").Value = '  +2.80%  '

# Row 6
$ws.Range("D6").Value = '104.92'
$ws.Range("E6").Value = '  +4.14%  '

# Row 7
$ws.Range("E7").Value = '  +1.17%  '

# Row 8
$ws.Range("E8").Value = '  -0.07%  '

# Row 9
$ws.Range("E9").Value = '  +5.03%  '

# Row 10
$ws.Range("D10").Value = '36.01'
$ws.Range("E10").Value = '  +2.50%  '

# Row 11
$ws.Range("D11").Value = '0.0807'
$ws.Range("E11").Value = '  +1.15%  '

# Row 12
$ws.Range("E12").Value = '  -1.60%  '

# Row 13
$ws.Range("D13").Value = '18.33'
$ws.Range("E13").Value = '  -4.03%  '

# Row 14
$ws.Range("D14").Value = '7.07'
$ws.Range("E14").Value = '  +2.28%  '

# Row 15
$ws.Range("D15").Value = '2.833.42'
$ws.Range("E15").Value = '  +1.27%  '

# Row 16
$ws.Range("D16").Value = '2.487.79'
$ws.Range("E16").Value = '  +2.28%  '

# Row 17
$ws.Range("D17").Value = '0.843'
$ws.Range("E17").Value = '  +1.45%  '

# Row 18
$ws.Range("D18").Value = '46.075.57'
$ws.Range("E18").Value = '  +4.03%  '

# Row 19
$ws.Range("D19").Value = '12.68'
$ws.Range("E19").Value = '  +2.38%  '

# Row 20
$ws.Range("D20").Value = '6.42'
$ws.Range("E20").Value = '  +0.97%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0935'
$ws.Range("E21").Value = '  +1.75%  '

# Row 22
$ws.Range("D22").Value = '70.94'
$ws.Range("E22").Value = '  +3.43%  '

# Row 23
$ws.Range("D23").Value = '2.40'
$ws.Range("E23").Value = '  +5.40%  '

# Row 24
$ws.Range("D24").Value = '247.39'
$ws.Range("E24").Value = '  +2.44%  '

# Row 25
$ws.Range("D25").Value = '2.52'
$ws.Range("E25").Value = '  +2.38%  '

# Row 26
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '25.92'
$ws.Range("E26").Value = '  +3.36%  '

# Row 27
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.05%  '

# Row 28
$ws.Range("D28").Value = '2.19'
$ws.Range("E28").Value = '  -3.64%  '

# Row 29
$ws.Range("E29").Value = '  +1.58%  '

# Row 30
$ws.Range("D30").Value = '34.47'
$ws.Range("E30").Value = '  +3.98%  '

# Row 31
$ws.Range("D31").Value = '49.41'
$ws.Range("E31").Value = '  +1.87%  '

# Row 32
$ws.Range("D32").Value = '0.128'
$ws.Range("E32").Value = '  +3.92%  '

# Row 33
$ws.Range("E33").Value = '  +2.80%  '

# Row 34
$ws.Range("D34").Value = '5.34'
$ws.Range("E34").Value = '  +3.65%  '

# Row 35
$ws.Range("E35").Value = '  -0.12%  '

# Row 36
$ws.Range("D36").Value = '0.0763'
$ws.Range("E36").Value = '  -0.98%  '

# Row 37
$ws.Range("D37").Value = '4.54'
$ws.Range("E37").Value = '  +1.28%  '

# Row 38
$ws.Range("E38").Value = '  +0.79%  '

# Row 39
$ws.Range("D39").Value = '2.96'
$ws.Range("E39").Value = '  +3.10%  '

# Row 40
$ws.Range("D40").Value = '127.32'
$ws.Range("E40").Value = '  +3.88%  '

# Row 41
$ws.Range("E41").Value = '  +1.90%  '

# Row 42
$ws.Range("D42").Value = '2.23'
$ws.Range("E42").Value = '  +0.48%  '

# Row 43
$ws.Range("D43").Value = '20.93'
$ws.Range("E43").Value = '  -1.11%  '

# Row 44
$ws.Range("D44").Value = '0.0293'
$ws.Range("E44").Value = '  +1.74%  '

# Row 45
$ws.Range("D45").Value = '1.971.62'
$ws.Range("E45").Value = '  +1.31%  '

# Row 46
$ws.Range("E46").Value = '  +1.82%  '

# Row 47
$ws.Range("D47").Value = '2.09'
$ws.Range("E47").Value = '  -4.23%  '

# Row 48
$ws.Range("E48").Value = '  +12.73%  '

# Row 49
$ws.Range("E49").Value = '  -4.07%  '

# Row 50
$ws.Range("D50").Value = '5.05'
$ws.Range("E50").Value = '  +8.68%  '

# Row 51
$ws.Range("D51").Value = '78.34'
$ws.Range("E51").Value = '  +5.94%  '

# Restore default style on column D so no stray formatting remains
$ws.Range("D2:D51").Style = "Normal"
